# Rows 16 and 17 in the sheet need to have their data swapped (the observation
# records that used to live on row 16 now belong on row 17, and vice versa).
# Only the columns whose content actually differs between the two rows are
# touched; columns that already hold identical values on both rows (D, I, T,
# U, V, W, Y, AA, AD, AE, AG, AT, AY) are left completely alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 16
$row2 = 17

# Columns whose value differs between row 16 and row 17 and therefore need
# to be exchanged.
$cols = @("A","B","E","F","G","H","K","L","M","N","P","Q","R","S","AC","AW","AX")

foreach ($col in $cols) {
    $addr1 = "$col$row1"
    $addr2 = "$col$row2"

    $val1 = $ws.Range($addr1).Value()
    $val2 = $ws.Range($addr2).Value()

    if ($val2 -eq $null) {
        $ws.Range($addr1).ClearContents()
    } else {
        $ws.Range($addr1).Value = $val2
    }

    if ($val1 -eq $null) {
        $ws.Range($addr2).ClearContents()
    } else {
        $ws.Range($addr2).Value = $val1
    }
}
